# Weekly fruit/vegetable price update: insert a new weekly record above the
# most recent "Crespo record"/Copenhague "Repollo" entries for Feria
# Lagunitas de Puerto Montt. This pushes the existing rows 579-613 down to
# 580-614 and fills the freed row 579 with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 579, shifting rows 579:613 down to 580:614.
$ws.Rows(579).Insert()

# Populate the newly inserted row 579 with this week's record.
$ws.Cells.Item(579, 1).Value = 4
$ws.Cells.Item(579, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(579, 3).Value = "Los Lagos"
$ws.Cells.Item(579, 4).Value = 44931
$ws.Cells.Item(579, 5).Value = 10
$ws.Cells.Item(579, 6).Value = 100112006
$ws.Cells.Item(579, 7).Value = "Repollo"
$ws.Cells.Item(579, 8).Value = "Crespo record"
$ws.Cells.Item(579, 9).Value = "Primera"
$ws.Cells.Item(579, 10).Value = 500
$ws.Cells.Item(579, 11).Value = 1700
$ws.Cells.Item(579, 12).Value = 1700
$ws.Cells.Item(579, 13).Value = 1700
$ws.Cells.Item(579, 14).Value = "$/unidad"
$ws.Cells.Item(579, 15).Value = "Región Metropolitana"
$ws.Cells.Item(579, 16).Value = 1700
$ws.Cells.Item(579, 17).Value = 1
$ws.Cells.Item(579, 18).Value = "Hortaliza"
